# Update "想去人数" (number of people interested) counts to the latest
# scraped values, per gh-pages data regeneration at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsShow    = $wb.Worksheets.Item(2)   # 演出
$wsLocal   = $wb.Worksheets.Item(3)   # 本地生活
$wsAll     = $wb.Worksheets.Item(4)   # 全部类型

# 展览 (sheet1)
$wsExhibit.Range("F5").Value  = 1123
$wsExhibit.Range("F8").Value  = 898
$wsExhibit.Range("F9").Value  = 1611
$wsExhibit.Range("F12").Value = 1741
$wsExhibit.Range("F14").Value = 5910
$wsExhibit.Range("F16").Value = 49
$wsExhibit.Range("F23").Value = 1373
$wsExhibit.Range("F24").Value = 720
$wsExhibit.Range("F25").Value = 238
$wsExhibit.Range("F28").Value = 27

# 演出 (sheet2)
$wsShow.Range("F4").Value = 310
$wsShow.Range("F5").Value = 160

# 本地生活 (sheet3)
$wsLocal.Range("F3").Value = 2229

# 全部类型 (sheet4)
$wsAll.Range("F3").Value  = 2229
$wsAll.Range("F7").Value  = 1123
$wsAll.Range("F11").Value = 310
$wsAll.Range("F12").Value = 898
$wsAll.Range("F14").Value = 1611
$wsAll.Range("F17").Value = 1741
$wsAll.Range("F23").Value = 5910
$wsAll.Range("F25").Value = 49
$wsAll.Range("F32").Value = 1373
$wsAll.Range("F33").Value = 720
$wsAll.Range("F35").Value = 238
